$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Published Values")
$ws2 = $wb.Worksheets.Item("Daily Summary")

# --- Append new row 6 to "All Published Values" ---
# Columns A, C, D, E, F, G hold date-like / numeric-like text that Excel would
# otherwise auto-convert to a date serial / number, so force them to Text
# first, then restore the default ("Normal") cell style so no stray number
# format sticks around on the new row (columns B, H, I, J round-trip fine as
# literal text without any extra help).

$r = $ws.Range("A6")
$r.NumberFormat = "@"
$r.Value = "2026-01-02"
$r.Style = "Normal"

$ws.Range("B6").Value = "2026-01-02 19:08:08"

$r = $ws.Range("C6")
$r.NumberFormat = "@"
$r.Value = "697.85"
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "697.85"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "700.79"
$r.Style = "Normal"

$r = $ws.Range("F6")
$r.NumberFormat = "@"
$r.Value = "700.79"
$r.Style = "Normal"

$r = $ws.Range("G6")
$r.NumberFormat = "@"
$r.Value = "702.88"
$r.Style = "Normal"

$ws.Range("H6").Value = "2026/01/02 19:08:08"
$ws.Range("I6").Value = "2026-01-02 11:10:50"
$ws.Range("J6").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# --- Re-establish the AutoFilter over the grown range A1:J6 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:J6").AutoFilter() | Out-Null

# --- Extend the hidden _FilterDatabase defined name for this sheet ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$6"
    }
}

# --- "Daily Summary": publishes count for 2026-01-02 goes from 4 to 5 ---
$ws2.Range("B4").Value = 5
